$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at row 12 to make room for "Docentes responsaveis:" + 4 names
# (this shifts "Programa resumido:" and everything below it down by 5 rows,
#  and the new rows inherit A/B/C column styles automatically)
$ws.Range("A12:A16").EntireRow.Insert()

# Fix "Objetivos:" / "Objectives:" body text (previously held the wrong content)
$ws.Range("B10").Value = 'Introdução às funções de variáveis complexas e suas aplicações. Apresentar equações diferenciais de interesse em engenharia física e desenvolver técnicas de soluções, verificando propriedades e métodos de resolução. Estudo de funções especiais em Engenharia Física.'
$ws.Range("C10").Value = 'Introdução às funções de variáveis complexas e suas aplicações. Apresentar equações diferenciais de interesse em engenharia física e desenvolver técnicas de soluções, verificando propriedades e métodos de resolução. Estudo de funções especiais em Engenharia Física.'

# New "Docentes responsaveis:" section with the 4 faculty names
$ws.Range("A12").Value = 'Docentes responsáveis:'
$ws.Range("B13").Value = '5840726 - Cristina Bormio Nunes'
$ws.Range("C13").Value = '5840726 - Cristina Bormio Nunes'
$ws.Range("B14").Value = '6495737 - Durval Rodrigues Junior'
$ws.Range("C14").Value = '6495737 - Durval Rodrigues Junior'
$ws.Range("B15").Value = '1341653 - Maria José Ramos Sandim'
$ws.Range("C15").Value = '1341653 - Maria José Ramos Sandim'
$ws.Range("B16").Value = '1643715 - Paulo Atsushi Suzuki'
$ws.Range("C16").Value = '1643715 - Paulo Atsushi Suzuki'

# "Programa resumido:" / "Short syllabus:" body text
$ws.Range("B17").Value = 'Funções de uma variável complexa. Função delta. Equações diferenciais parciais da engenharia física: métodos de solução, resolução de problemas de valores de contorno, aplicações. Série de Fourier e Transformadas Integrais. Funções especiais.'
$ws.Range("C17").Value = 'Funções de uma variável complexa. Função delta. Equações diferenciais parciais da engenharia física: métodos de solução, resolução de problemas de valores de contorno, aplicações. Série de Fourier e Transformadas Integrais. Funções especiais.'

# "Programa:" body text (Syllabus: body was already correct)
$ws.Range("B19").Value = 'Funções de uma variável complexa: séries infinitas, funções analíticas, condições de Cauchy Riemann, integrais de contorno, teorema de Cauchy, teorema dos resíduos, Função delta. Equação de Laplace, equação da difusão (do calor), equação de ondas (corda vibrante); Série de Fourier, Transformadas Integrais de Fourier e Laplace. Funções especiais: Polinômios de Legendre, Harmônicos Esféricos, Funções de Bessel.'
$ws.Range("C19").Value = 'Funções de uma variável complexa: séries infinitas, funções analíticas, condições de Cauchy Riemann, integrais de contorno, teorema de Cauchy, teorema dos resíduos, Função delta. Equação de Laplace, equação da difusão (do calor), equação de ondas (corda vibrante); Série de Fourier, Transformadas Integrais de Fourier e Laplace. Funções especiais: Polinômios de Legendre, Harmônicos Esféricos, Funções de Bessel.'

# "Metodo:", "Criterio:", "Norma de recuperacao:", "Bibliografia:" body text
$ws.Range("B22").Value = 'Aulas expositivas teóricas, aulas de exercícios.'
$ws.Range("C22").Value = 'Aulas expositivas teóricas, aulas de exercícios.'
$ws.Range("B23").Value = 'Duas provas escritas: conceitos P1 e P2. Conceito Final = (P1 + P2)/2'
$ws.Range("C23").Value = 'Duas provas escritas: conceitos P1 e P2. Conceito Final = (P1 + P2)/2'
$ws.Range("B24").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C24").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("B25").Value = '•ARFKEN, G. and WEBER, H. J. Mathematical Methods for Physicists.•BROWN, JAMES W. and CHURCHILL, RUEL V., Complex Variables and Applications, Mc Graw Hill Higher Education, 7a. ed.• BUTKOV, Eugene. Física Matemática.•BELLANDI FILHO,J., Funções Especiais, Ed. Papirus, 1985.'
$ws.Range("C25").Value = '•ARFKEN, G. and WEBER, H. J. Mathematical Methods for Physicists.•BROWN, JAMES W. and CHURCHILL, RUEL V., Complex Variables and Applications, Mc Graw Hill Higher Education, 7a. ed.• BUTKOV, Eugene. Física Matemática.•BELLANDI FILHO,J., Funções Especiais, Ed. Papirus, 1985.'
